$wb = $excel.ActiveWorkbook

# 1) Shared string change: "Ready for handoff" -> "In Translation"
#    (used by the Status-like columns on all three sheets: E/F on "Overview",
#    C on "zh-cn" and "de-de")
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2) Narrow the "status" column(s) from ~17.22 chars to ~13.41 chars.
#    ColumnWidth is quantized internally to 1/6-character steps, so 12.5 is
#    the input that lands closest to the target stored width (13.41...).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)
